# "Updated solution to 125" - adds a new row for LeetCode problem 27
# ("Remove Element"), and widens/restyles the Tags & Solution columns
# (left-aligned, vertically centered) to accommodate the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new row (problem 27: "remove element") -----------------------
$ws.Range("A4").Value = 27
$ws.Range("B4").Value = "remove element"
$ws.Range("C4").Value = "Easy"
$ws.Range("D4").Value = "Two pointers; Sequential"
$ws.Range("E4").Value = "Initialize two pointers one at the start and one at the end, replace left element with right element if equals to the value"

# --- Re-align the Tags (D) / Solution (E) columns --------------------------
# Apply vertical-center first, then horizontal-left, so the resulting style
# keeps the vertical centering that was already present on some of these
# cells. Data rows (normal font) are done before the header row (bold font)
# so new style entries land in the same order as authored in Excel.
$ws.Range("D3:E4").VerticalAlignment = -4108   # xlCenter
$ws.Range("D3:E4").HorizontalAlignment = -4131 # xlLeft

$ws.Range("D1:E1").VerticalAlignment = -4108   # xlCenter
$ws.Range("D1:E1").HorizontalAlignment = -4131 # xlLeft

# --- Widen the Solution column to fit the longer text ----------------------
$ws.Columns("E").ColumnWidth = 94.875

# --- Move the active cell selection down to the new next-empty row --------
$ws.Range("A5").Select() | Out-Null

Write-Host "Added problem 27 (remove element) row and updated column formatting"
